$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# G118 before: style 2 (no alignment override); want style 1 (right aligned) -- copy from G119/G120 which already have style=1
$ws.Range("G120").Copy()
$ws.Range("G118").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# G119 before: style 1 (right aligned); want style 2 (no override) -- copy from G118/G121 which already have style=2
$ws.Range("G121").Copy()
$ws.Range("G119").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
